$wb = $excel.ActiveWorkbook

# Grab the two worksheets by their current names (order-independent).
$hotelWs  = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")

# --- Update hotel_info: insert a new "State" column right after "Hotel_Name" ---
# Hotel_Name is column B, so the new column goes in at column C and everything
# from City onward shifts one column to the right.
$hotelWs.Columns.Item(3).Insert()
$hotelWs.Range("C1").Value = "State"
$hotelWs.Range("C2").Value = "Louisiana"

# --- Reorder the sheets: review_info becomes the first sheet, hotel_info second ---
$reviewWs.Move($hotelWs, $null)
